# Record today's progress (04-Nov-2025): for each training row on the
# "Training Dashboard" sheet, the "PERIOD TO EXPIRE" (column H) drops by one
# day and the "LAST UPDATE" (column I) moves from 03-Nov-2025 to 04-Nov-2025.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Training Dashboard")

$newLastUpdate = "04-Nov-2025"
$firstRow = 3
$lastRow = 14

# Force column I to stay text so the date string is not reinterpreted as a
# serial date value when it is written below.
$ws.Range("I$firstRow`:I$lastRow").NumberFormat = "@"

for ($row = $firstRow; $row -le $lastRow; $row++) {
    $hCell = $ws.Cells.Item($row, 8)  # PERIOD TO EXPIRE
    $hCell.Value2 = $hCell.Value2 - 1

    $iCell = $ws.Cells.Item($row, 9)  # LAST UPDATE
    $iCell.Value2 = $newLastUpdate
}
